$d = $word.ActiveDocument

# --- Change 1: "Low correlation between uncertainties" -> "Low correlation between model uncertainties"
$found1 = $d.Content.Find.Execute(
    "Low correlation between uncertainties estimated using the two methods",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Low correlation between model uncertainties estimated using the two methods",
    2)

# --- Change 2: remove " for the four model parameters" before "shows they are not equivalent"
$found2 = $d.Content.Find.Execute(
    " for the four model parameters shows they are not equivalent",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " shows they are not equivalent",
    2)

# --- Relocate the hidden "_GoBack" bookmark to sit inside "those" ("th" | "ose"),
#     matching where Word leaves it after the most recent edit in this sentence.
$text = $d.Content.Text
$idx = $text.IndexOf("a wider range than those estimated simply")
if ($idx -ge 0) {
    $splitPos = $idx + "a wider range than th".Length
    $bmRange = $d.Range($splitPos, $splitPos)
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}

Write-Output "found1=$found1 found2=$found2"
